$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 5, 6, 7 duplicate rows 4, 2, 3 respectively (same venue/date/
# result/teams/batsman/figures), extending the table the same way the
# source scrape re-appended match rows.
$newRows = @(
    @(" Dubai (DSC)", " October 13 2020", "Super Kings won by 20 runs", "Chennai Super Kings", "Sunrisers Hyderabad", "Deepak Chahar ", "2", "2", "0", "0", "100.00"),
    @(" Sharjah", " October 23 2020", "Mumbai won by 10 wickets (with 46 balls remaining)", "Chennai Super Kings", "Mumbai Indians", "Deepak Chahar ", "0", "5", "0", "0", "0.00"),
    @(" Dubai (DSC)", " October 10 2020", "RCB won by 37 runs", "Chennai Super Kings", "Royal Challengers Bangalore", "Deepak Chahar ", "5", "5", "0", "0", "100.00")
)

$startRow = 5
$endRow = $startRow + $newRows.Count - 1
$numCols = $newRows[0].Count

# Every existing cell in the sheet (rows 1-4) is stored as text, even the
# numeric-looking figures in columns G:K ("0", "2", "100.00", ...) - that's
# why the sheet also carries a numberStoredAsText ignoredError over the
# whole table. Plain `.Value = "2"` assignment auto-converts such strings
# to real numbers, so the target range is put into Text format first, the
# literal values are written (so they stay text), and then the format is
# reset back to the workbook default so the new cells don't end up with a
# different style than the rest of the table.
$targetRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, $numCols))
$targetRange.NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

$targetRange.ClearFormats()
